# Update the 'fo high low' data table (Sheet1) with new High/Low/Close/LTP/Vol/9:25 Close values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 797.6
$ws.Range("C2").Value = 773.6
$ws.Range("D2").Value = 787
$ws.Range("E2").Value = 790.35
$ws.Range("F2").Value = 192
$ws.Range("G2").Value = 786.7

# Row 3
$ws.Range("B3").Value = 889.6
$ws.Range("C3").Value = 856.15
$ws.Range("D3").Value = 876
$ws.Range("E3").Value = 875.5
$ws.Range("F3").Value = 91
$ws.Range("G3").Value = 860.7

# Row 4
$ws.Range("B4").Value = 44347.55
$ws.Range("C4").Value = 43982.4
$ws.Range("D4").Value = 44219.7
$ws.Range("E4").Value = 44234.05
$ws.Range("F4").Value = 26
$ws.Range("G4").Value = 44025

# Row 5
$ws.Range("B5").Value = 332.7
$ws.Range("C5").Value = 323.65
$ws.Range("D5").Value = 329.55
$ws.Range("E5").Value = 330.8
$ws.Range("F5").Value = 192
$ws.Range("G5").Value = 327.5

# Row 6
$ws.Range("B6").Value = 474.8
$ws.Range("C6").Value = 464.8
$ws.Range("D6").Value = 472
$ws.Range("E6").Value = 472.95
$ws.Range("F6").Value = 103
$ws.Range("G6").Value = 471.4

# Row 7
$ws.Range("B7").Value = 454.6
$ws.Range("C7").Value = 448.6
$ws.Range("D7").Value = 451.1
$ws.Range("E7").Value = 451.6
$ws.Range("F7").Value = 71
$ws.Range("G7").Value = 449.85

# Row 8
$ws.Range("B8").Value = 964.15
$ws.Range("C8").Value = 951.45
$ws.Range("D8").Value = 960.9
$ws.Range("E8").Value = 961.2
$ws.Range("F8").Value = 122
$ws.Range("G8").Value = 953.2

# Row 9
$ws.Range("B9").Value = 680.45
$ws.Range("C9").Value = 657.3
$ws.Range("D9").Value = 662
$ws.Range("E9").Value = 662.4
$ws.Range("F9").Value = 123
$ws.Range("G9").Value = 677.7

# Row 10
$ws.Range("B10").Value = 19506.25
$ws.Range("C10").Value = 19316.85
$ws.Range("D10").Value = 19471.2
$ws.Range("E10").Value = 19478.35
$ws.Range("F10").Value = 57
$ws.Range("G10").Value = 19334.15

# Row 11
$ws.Range("B11").Value = 2581.75
$ws.Range("C11").Value = 2527.8
$ws.Range("D11").Value = 2576
$ws.Range("E11").Value = 2577.45
$ws.Range("F11").Value = 97
$ws.Range("G11").Value = 2529.7

# Row 12
$ws.Range("B12").Value = 567.1
$ws.Range("C12").Value = 560.75
$ws.Range("D12").Value = 563
$ws.Range("E12").Value = 563.2
$ws.Range("F12").Value = 314
$ws.Range("G12").Value = 565.9

# Row 13
$ws.Range("B13").Value = 847.1
$ws.Range("C13").Value = 832.6
$ws.Range("D13").Value = 846.5
$ws.Range("E13").Value = 845.95
$ws.Range("F13").Value = 14
$ws.Range("G13").Value = 832.85

# Row 14
$ws.Range("B14").Value = 610.15
$ws.Range("C14").Value = 596.25
$ws.Range("D14").Value = 607.95
$ws.Range("E14").Value = 608.25
$ws.Range("F14").Value = 196
$ws.Range("G14").Value = 598.55

# Row 15
$ws.Range("B15").Value = 119.15
$ws.Range("C15").Value = 117.5
$ws.Range("D15").Value = 118.35
$ws.Range("E15").Value = 118.45
$ws.Range("F15").Value = 419
$ws.Range("G15").Value = 118.95

# Row 16
$ws.Range("B16").Value = 3461.45
$ws.Range("C16").Value = 3423
$ws.Range("D16").Value = 3458.4
$ws.Range("E16").Value = 3457.45
$ws.Range("F16").Value = 14
$ws.Range("G16").Value = 3438.25

# Row 17
$ws.Range("B17").Value = 3026.75
$ws.Range("C17").Value = 2996.05
$ws.Range("D17").Value = 3023
$ws.Range("E17").Value = 3021.35
$ws.Range("F17").Value = 9
$ws.Range("G17").Value = 3008.05

# Update the active cell / selection to match the saved view state
[void]$ws.Range("K14").Select()
